# LH_TC_SYSTEMCONSTRAINTS_REVIEWS.xlsx
# v2.1 - closed the review comments on v2.0

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Version History")
$ws2 = $wb.Worksheets.Item("LH-TC-SYSTEMCONSTRAINS-REVIEWS")

# --- Version History sheet -------------------------------------------------

# Update the existing v2.0 entry's "Updated Section" text
$ws1.Range("C5").Value = "Reviewed on v2.0 according to new srs"

# Grow the Table1 ListObject by one row (keeps the table definition's ref in sync)
$lo = $ws1.ListObjects.Item("Table1")
$lo.ListRows.Add() | Out-Null

# Copy formatting (styles + row height) from row 5 down into the new row 6
$ws1.Range("A5:D5").Copy()
$ws1.Range("A6:D6").PasteSpecial(-4122)
$ws1.Rows(6).RowHeight = 28.8

# New version-history entry for v2.1
$ws1.Range("A6").Value = "v2.1"
$ws1.Range("B6").Value = "Mahmoud Abdelmageed"
$ws1.Range("C6").Value = "closed the reviews on v2.0"
$ws1.Range("D6").Formula = "=DATE(2025,5,11)"

# Move the sheet1 selection (no longer the active tab)
$ws1.Range("B12").Select()

# --- LH-TC-SYSTEMCONSTRAINS-REVIEWS sheet -----------------------------------

# Widen column H (drop the old "best fit" auto width)
$ws2.Columns(8).ColumnWidth = 19.5

# Row 4 height shrinks slightly
$ws2.Rows(4).RowHeight = 189

# This sheet becomes the active / selected tab
$ws2.Activate()
